$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 19
$ws.Range("A19").Value = 1700
$ws.Range("C18").Copy()
$ws.Range("C19").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C19").Value = 44774
$ws.Range("D19").Value = "UVC LAMP"

# Row 20
$ws.Range("A20").Value = 202
$ws.Range("D20").Value = "petri dish"

# Row 21
$ws.Range("A21").Value = 257
$ws.Range("D21").Value = "premixed nutrient agar"

$ws.Range("D21").Select()
$excel.ActiveWindow.Zoom = 130
